$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.372729943998593
$ws.Range("C2").Value = 0.2344453039142707
$ws.Range("D2").Value = 0.1496931498150786
$ws.Range("F2").Value = 1.747409726797848
$ws.Range("G2").Value = 0.002491246750073706
$ws.Range("J2").Value = 0.2184199188407483
$ws.Range("L2").Value = 0.3277370905224615
$ws.Range("M2").Value = 0.3342352116041098
$ws.Range("O2").Value = 4.440290433129405
$ws.Range("B3").Value = 1.273274568625254
$ws.Range("C3").Value = 0.222118127250809
$ws.Range("D3").Value = 0.1495340626660635
$ws.Range("F3").Value = 1.759505153369659
$ws.Range("G3").Value = 0.002494296656742231
$ws.Range("J3").Value = 0.2202560094990664
$ws.Range("L3").Value = 0.3243010243868767
$ws.Range("M3").Value = 0.3183151109437574
$ws.Range("O3").Value = 4.475366979812634
$ws.Range("B4").Value = 1.212401018146977
$ws.Range("C4").Value = 0.2145167690970311
$ws.Range("D4").Value = 0.1494715267702489
$ws.Range("F4").Value = 1.767931329572093
$ws.Range("G4").Value = 0.002496269968750262
$ws.Range("J4").Value = 0.2214464610683331
$ws.Range("L4").Value = 0.3222924456805032
$ws.Range("M4").Value = 0.3086180383821429
$ws.Range("O4").Value = 4.499649157197808
$ws.Range("B5").Value = 1.187644517254512
$ws.Range("C5").Value = 0.2114112111658812
$ws.Range("D5").Value = 0.1494549141167205
$ws.Range("F5").Value = 1.771616364344361
$ws.Range("G5").Value = 0.002497099496608449
$ws.Range("J5").Value = 0.2219474673145161
$ws.Range("L5").Value = 0.3214994902199706
$ws.Range("M5").Value = 0.3046862777406858
$ws.Range("O5").Value = 4.510234191781194
$ws.Range("B6").Value = 1.183536783455565
$ws.Range("C6").Value = 0.2108950621402954
$ws.Range("D6").Value = 0.1494526923988175
$ws.Range("F6").Value = 1.772243438218432
$ws.Range("G6").Value = 0.00249723877459294
$ws.Range("J6").Value = 0.2220316193118013
$ws.Range("L6").Value = 0.3213693675025979
$ws.Range("M6").Value = 0.304034620523602
$ws.Range("O6").Value = 4.512033481905405
$ws.Range("B7").Value = 1.212066938917815
$ws.Range("C7").Value = 0.2144749183346732
$ws.Range("D7").Value = 0.1494712667619815
$ws.Range("F7").Value = 1.767980009796098
$ws.Range("G7").Value = 0.002496281053172756
$ws.Range("J7").Value = 0.2214531534531403
$ws.Range("L7").Value = 0.3222816479755934
$ws.Range("M7").Value = 0.3085649324389053
$ws.Range("O7").Value = 4.499789118055645
$ws.Range("B8").Value = 1.338398694043633
$ws.Range("C8").Value = 0.2302017441235478
$ws.Range("D8").Value = 0.1496310224307393
$ws.Range("F8").Value = 1.75137274988225
$ws.Range("G8").Value = 0.002492277512373658
$ws.Range("J8").Value = 0.2190399245280119
$ws.Range("L8").Value = 0.3265314017495058
$ws.Range("M8").Value = 0.3287299466218698
$ws.Range("O8").Value = 4.451814844256717
$ws.Range("B9").Value = 1.587606049736621
$ws.Range("C9").Value = 0.260777055885427
$ws.Range("D9").Value = 0.1502217721476242
$ws.Range("F9").Value = 1.726739608662385
$ws.Range("G9").Value = 0.002485221773770644
$ws.Range("J9").Value = 0.2148069951574767
$ws.Range("L9").Value = 0.3356636365308532
$ws.Range("M9").Value = 0.3688818483715295
$ws.Range("O9").Value = 4.379538598096104
$ws.Range("B10").Value = 1.771539926112723
$ws.Range("C10").Value = 0.2830711338289404
$ws.Range("D10").Value = 0.1508232430850782
$ws.Range("F10").Value = 1.713483357448993
$ws.Range("G10").Value = 0.002480517821750525
$ws.Range("J10").Value = 0.2119998781291788
$ws.Range("L10").Value = 0.342855040756362
$ws.Range("M10").Value = 0.3987421647247302
$ws.Range("O10").Value = 4.339759101514147
$ws.Range("B11").Value = 1.855388475378504
$ws.Range("C11").Value = 0.2931749188502693
$ws.Range("D11").Value = 0.1511328785201513
$ws.Range("F11").Value = 1.708505290158797
$ws.Range("G11").Value = 0.002478481040890523
$ws.Range("J11").Value = 0.2107882549163982
$ws.Range("L11").Value = 0.3462302986952324
$ws.Range("M11").Value = 0.4124028010190486
$ws.Range("O11").Value = 4.324562051294009
$ws.Range("B12").Value = 1.887163786622295
$ws.Range("C12").Value = 0.2969953348789716
$ws.Range("D12").Value = 0.1512552779960714
$ws.Range("F12").Value = 1.706771624292429
$ws.Range("G12").Value = 0.00247772450828481
$ws.Range("J12").Value = 0.210338817236071
$ws.Range("L12").Value = 0.3475232561057737
$ws.Range("M12").Value = 0.417586576630022
$ws.Range("O12").Value = 4.319224710791957
$ws.Range("B13").Value = 1.880319377335127
$ws.Range("C13").Value = 0.2961727954588866
$ws.Range("D13").Value = 0.151228688674145
$ws.Range("F13").Value = 1.707138264132652
$ws.Range("G13").Value = 0.00247788678597937
$ws.Range("J13").Value = 0.2104351949085501
$ws.Range("L13").Value = 0.347244137245454
$ws.Range("M13").Value = 0.4164696825065661
$ws.Range("O13").Value = 4.320355626676957
$ws.Range("B14").Value = 1.858002184002657
$ws.Range("C14").Value = 0.293489341612343
$ws.Range("D14").Value = 0.1511428454272306
$ws.Range("F14").Value = 1.708359625138016
$ws.Range("G14").Value = 0.002478418505344557
$ws.Range("J14").Value = 0.2107510915735666
$ws.Range("L14").Value = 0.346336374701238
$ws.Range("M14").Value = 0.4128290587668886
$ws.Range("O14").Value = 4.324114575263735
$ws.Range("B15").Value = 1.844335294842949
$ws.Range("C15").Value = 0.2918449040770383
$ws.Range("D15").Value = 0.1510909333137107
$ws.Range("F15").Value = 1.709127466510566
$ws.Range("G15").Value = 0.002478746117765584
$ws.Range("J15").Value = 0.2109458080929816
$ws.Range("L15").Value = 0.345782270645401
$ws.Range("M15").Value = 0.410600468174465
$ws.Range("O15").Value = 4.326471422305673
$ws.Range("B16").Value = 1.766063629544362
$ws.Range("C16").Value = 0.2824100463463708
$ws.Range("D16").Value = 0.1508037297490858
$ws.Range("F16").Value = 1.71382986919447
$ws.Range("G16").Value = 0.002480652998202043
$ws.Range("J16").Value = 0.2120803738765353
$ws.Range("L16").Value = 0.3426365397779989
$ws.Range("M16").Value = 0.3978509353991768
$ws.Range("O16").Value = 4.3408106481684
$ws.Range("B17").Value = 1.718090432796259
$ws.Range("C17").Value = 0.2766122031780185
$ws.Range("D17").Value = 0.1506367426556423
$ws.Range("F17").Value = 1.716984227159998
$ws.Range("G17").Value = 0.002481849156705245
$ws.Range("J17").Value = 0.2127931167125254
$ws.Range("L17").Value = 0.3407332539312193
$ws.Range("M17").Value = 0.3900490335728009
$ws.Range("O17").Value = 4.350350191355943
$ws.Range("B18").Value = 1.690514191861837
$ws.Range("C18").Value = 0.273273882387798
$ws.Range("D18").Value = 0.1505440887307401
$ws.Range("F18").Value = 1.718897568603751
$ws.Range("G18").Value = 0.002482546861223401
$ws.Range("J18").Value = 0.2132092195739208
$ws.Range("L18").Value = 0.3396483160322958
$ws.Range("M18").Value = 0.385568853330561
$ws.Range("O18").Value = 4.356109901133465
$ws.Range("B19").Value = 1.681180264271177
$ws.Range("C19").Value = 0.2721429804726085
$ws.Range("D19").Value = 0.1505133013960247
$ws.Range("F19").Value = 1.719562399826046
$ws.Range("G19").Value = 0.002482784760796944
$ws.Range("J19").Value = 0.2133511620638613
$ws.Range("L19").Value = 0.3392826581961401
$ws.Range("M19").Value = 0.3840531982509745
$ws.Range("O19").Value = 4.358106878072249
$ws.Range("B20").Value = 1.723195547809837
$ws.Range("C20").Value = 0.277229762968858
$ws.Range("D20").Value = 0.1506541678311706
$ws.Range("F20").Value = 1.716638189521987
$ws.Range("G20").Value = 0.002481720819573053
$ws.Range("J20").Value = 0.2127166074486944
$ws.Range("L20").Value = 0.3409348504053327
$ws.Range("M20").Value = 0.3908788092893545
$ws.Range("O20").Value = 4.349306450795467
$ws.Range("B21").Value = 1.864556654833791
$ws.Range("C21").Value = 0.2942776926067552
$ws.Range("D21").Value = 0.1511679202360199
$ws.Range("F21").Value = 1.707996771468189
$ws.Range("G21").Value = 0.002478261926568393
$ws.Range("J21").Value = 0.210658050708874
$ws.Range("L21").Value = 0.3466026052412019
$ws.Range("M21").Value = 0.413898107675351
$ws.Range("O21").Value = 4.322999146610329
$ws.Range("B22").Value = 1.957081616352582
$ws.Range("C22").Value = 0.3053863630255194
$ws.Range("D22").Value = 0.1515336741803353
$ws.Range("F22").Value = 1.7032317463064
$ws.Range("G22").Value = 0.002476087293089967
$ws.Range("J22").Value = 0.2093673167541414
$ws.Range("L22").Value = 0.3503931437022345
$ws.Range("M22").Value = 0.4290053053451217
$ws.Range("O22").Value = 4.308239313848503
$ws.Range("B23").Value = 1.907687277219964
$ws.Range("C23").Value = 0.2994605598666737
$ws.Range("D23").Value = 0.1513357313996977
$ws.Range("F23").Value = 1.705694133574042
$ws.Range("G23").Value = 0.002477240094025782
$ws.Range("J23").Value = 0.2100512114016713
$ws.Range("L23").Value = 0.3483621997638835
$ws.Range("M23").Value = 0.4209366618931298
$ws.Range("O23").Value = 4.315894049413828
$ws.Range("B24").Value = 1.72088751214153
$ws.Range("C24").Value = 0.2769505799444687
$ws.Range("D24").Value = 0.1506462794748771
$ws.Range("F24").Value = 1.716794322087736
$ws.Range("G24").Value = 0.002481778809622058
$ws.Range("J24").Value = 0.2127511775438027
$ws.Range("L24").Value = 0.340843679705003
$ws.Range("M24").Value = 0.390503651353832
$ws.Range("O24").Value = 4.349777468394649
$ws.Range("B25").Value = 1.5200371660992
$ws.Range("C25").Value = 0.2525348341553411
$ws.Range("D25").Value = 0.1500324307933809
$ws.Range("F25").Value = 1.732553639827103
$ws.Range("G25").Value = 0.002487045912782869
$ws.Range("J25").Value = 0.2158988200071506
$ws.Range("L25").Value = 0.3331081437341226
$ws.Range("M25").Value = 0.3579556248875164
$ws.Range("O25").Value = 4.39675403231638
